$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to text format first so numeric-looking strings
# like "62.920.30" are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.920.30'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').Value = '3.408.47'
$ws.Range('E3').Value = '  -3.01%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '575.34'
$ws.Range('E5').Value = '  -2.81%  '
$ws.Range('D6').Value = '126.43'
$ws.Range('E6').Value = '  -6.11%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.407.91'
$ws.Range('E8').Value = '  -3.01%  '
$ws.Range('E9').Value = '  -2.59%  '
$ws.Range('D10').Value = '7.38'
$ws.Range('E10').Value = '  -3.22%  '
$ws.Range('E11').Value = '  -2.90%  '
$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -2.91%  '
$ws.Range('D13').Value = '3.980.86'
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').Value = '3.400.37'
$ws.Range('E15').Value = '  -3.24%  '
$ws.Range('E16').Value = '  -4.49%  '
$ws.Range('D17').Value = '62.902.65'
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').Value = '24.86'
$ws.Range('E18').Value = '  -3.86%  '
$ws.Range('D19').Value = '9.55'
$ws.Range('E19').Value = '  -4.47%  '
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').Value = '13.18'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Value = '376.88'
$ws.Range('E22').Value = '  -4.37%  '
$ws.Range('D23').Value = '0.559'
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').Value = '3.539.44'
$ws.Range('E24').Value = '  -3.21%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '71.96'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('E27').Value = '  -7.75%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = '7.00'
$ws.Range('E29').Value = '  -5.48%  '
$ws.Range('D30').Value = '2.16'
$ws.Range('E30').Value = '  -4.73%  '
$ws.Range('D31').Value = '7.86'
$ws.Range('E31').Value = '  -5.34%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.151'
$ws.Range('E32').Value = '  -3.50%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.40'
$ws.Range('E33').Value = '  -4.70%  '
$ws.Range('D35').Value = '3.433.30'
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('D36').Value = '22.78'
$ws.Range('E36').Value = '  -2.75%  '
$ws.Range('D37').Value = '5.30'
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('D38').Value = '165.28'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('E39').Value = '  -3.50%  '
$ws.Range('E40').Value = '  -4.20%  '
$ws.Range('D41').Value = '0.0760'
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '0.776'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('D44').Value = '41.63'
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('D45').Value = '4.29'
$ws.Range('E45').Value = '  -3.65%  '
$ws.Range('D46').Value = '1.58'
$ws.Range('E46').Value = '  -5.32%  '
$ws.Range('D47').Value = '22.91'
$ws.Range('E47').Value = '  -10.25%  '
$ws.Range('D48').Value = '1.09'
$ws.Range('E48').Value = '  -7.92%  '
$ws.Range('D49').Value = '6.67'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').Value = '2.254.20'
$ws.Range('E50').Value = '  -5.88%  '
$ws.Range('D51').Value = '0.857'
$ws.Range('E51').Value = '  -4.65%  '

# Restore the default (unstyled) look so no stray style survives the edit.
$ws.Range("D2:D51").Style = "Normal"
